# Update loading_percent values for the 380 kV case (rows 2-25, columns B,D-L,N)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    "B2" = 14.69890608103378
    "D2" = 7.572163081662111
    "E2" = 12.64670359564358
    "F2" = 38.65643293357294
    "G2" = 44.83942540982797
    "H2" = 18.28927540083934
    "I2" = 28.87748646452587
    "J2" = 10.14041089272122
    "K2" = 12.69257957617437
    "L2" = 11.26114052038591
    "N2" = 21.19674134173647
    "B3" = 14.59064597173884
    "D3" = 7.574029567663559
    "E3" = 12.66228181817337
    "F3" = 38.70234376265739
    "G3" = 44.83744534736566
    "H3" = 18.33097746650159
    "I3" = 28.97095124947406
    "J3" = 10.1507600937019
    "K3" = 12.4847955616569
    "L3" = 11.23446620043975
    "N3" = 21.26055695412704
    "B4" = 14.52645458149947
    "D4" = 7.575680571815174
    "E4" = 12.67252805241723
    "F4" = 38.73876635179204
    "G4" = 44.84819120245588
    "H4" = 18.35960780550905
    "I4" = 29.03264832855853
    "J4" = 10.15748384155383
    "K4" = 12.35816043131097
    "L4" = 11.21993807970401
    "N4" = 21.30155498930459
    "B5" = 14.50089345629246
    "D5" = 7.576480812249745
    "E5" = 12.67687514958782
    "F5" = 38.75567677409581
    "G5" = 44.85557615107563
    "H5" = 18.37203495458797
    "I5" = 29.05887404913103
    "J5" = 10.16031692367345
    "K5" = 12.30686123486642
    "L5" = 11.21448717109725
    "N5" = 21.31871977396381
    "B6" = 14.49668580996793
    "D6" = 7.576621402159354
    "E6" = 12.67760736206126
    "F6" = 38.7586095516034
    "G6" = 44.85698387117777
    "H6" = 18.37414435889886
    "I6" = 29.06329426212328
    "J6" = 10.1607929848339
    "K6" = 12.29836349182864
    "L6" = 11.21361051931387
    "N6" = 21.32159766044909
    "B7" = 14.52610740562174
    "D7" = 7.575690847507819
    "E7" = 12.67258598325073
    "F7" = 38.7389860424016
    "G7" = 44.84827863166728
    "H7" = 18.35977232586252
    "I7" = 29.03299763004532
    "J7" = 10.15752167220975
    "K7" = 12.35746726494149
    "L7" = 11.2198626609053
    "N7" = 21.30178462438078
    "B8" = 14.66111923705138
    "D8" = 7.572702118167387
    "E8" = 12.65193383255105
    "F8" = 38.67055296301885
    "G8" = 44.83625948317012
    "H8" = 18.30302604951833
    "I8" = 28.90881874987712
    "J8" = 10.14390279284826
    "K8" = 12.62078021055829
    "L8" = 11.25156188343869
    "N8" = 21.21836908510122
    "B9" = 14.94288199565255
    "D9" = 7.570827814412924
    "E9" = 12.61682201713771
    "F9" = 38.6017633963126
    "G9" = 44.90760140315315
    "H9" = 18.21577646473907
    "I9" = 28.69949479323208
    "J9" = 10.1201156168994
    "K9" = 13.14167016186625
    "L9" = 11.32820228472012
    "N9" = 21.06912947142988
    "B10" = 15.15874476450094
    "D10" = 7.571854270380575
    "E10" = 12.59428544982521
    "F10" = 38.59117156354175
    "G10" = 45.01775520206738
    "H10" = 18.16635784018441
    "I10" = 28.56654492897672
    "J10" = 10.10440387935867
    "K10" = 13.52323290688246
    "L10" = 11.39304525918605
    "N10" = 20.96813393469831
    "B11" = 15.2585338055182
    "D11" = 7.572837245370111
    "E11" = 12.58473583760695
    "F11" = 38.59502795105966
    "G11" = 45.08033418377342
    "H11" = 18.14707053885883
    "I11" = 28.51058769950385
    "J11" = 10.09763620632826
    "K11" = 13.69577544235407
    "L11" = 11.42432793955922
    "N11" = 20.92404765668586
    "B12" = 15.29652155568308
    "D12" = 7.573283149491148
    "E12" = 12.5812202605205
    "F12" = 38.59773425198453
    "G12" = 45.10581464712761
    "H12" = 18.14022652720936
    "I12" = 28.49004857545084
    "J12" = 10.09512782163849
    "K12" = 13.76090118216321
    "L12" = 11.43642415991235
    "N12" = 20.90761896001167
    "B13" = 15.28833179299126
    "D13" = 7.573183848053711
    "E13" = 12.58197293185603
    "F13" = 38.59709601836948
    "G13" = 45.10024784844637
    "H13" = 18.14168005727159
    "I13" = 28.49444309896799
    "J13" = 10.09566563165184
    "K13" = 13.74688562331965
    "L13" = 11.433808005102
    "N13" = 20.91114537081345
    "B14" = 15.26165524357631
    "D14" = 7.572872456749888
    "E14" = 12.58444459369389
    "F14" = 38.59522564010033
    "G14" = 45.08239480434973
    "H14" = 18.14649826312063
    "I14" = 28.508884891019
    "J14" = 10.0974287512275
    "K14" = 13.70113798902936
    "L14" = 11.42531812834013
    "N14" = 20.92269073749245
    "B15" = 15.24534024018536
    "D15" = 7.572691298953693
    "E15" = 12.58597165575824
    "F15" = 38.59424218494591
    "G15" = 45.07169118714278
    "H15" = 18.14950942803238
    "I15" = 28.51781565224132
    "J15" = 10.09851578935506
    "K15" = 13.67308673333055
    "L15" = 11.4201502097676
    "N15" = 20.92979718928618
    "B16" = 15.15225286609749
    "D16" = 7.571800369369459
    "E16" = 12.59492364337295
    "F16" = 38.5910939919801
    "G16" = 45.0139156372873
    "H16" = 18.16768261220705
    "I16" = 28.57029289669935
    "J16" = 10.10485378468648
    "K16" = 13.51193054841113
    "L16" = 11.39103623245704
    "N16" = 20.97105234094988
    "B17" = 15.09553387237149
    "D17" = 7.571385606962313
    "E17" = 12.600595046546
    "F17" = 38.59138370587686
    "G17" = 44.9816595725838
    "H17" = 18.17964954671744
    "I17" = 28.6036445143268
    "J17" = 10.10883903286888
    "K17" = 13.41275573501913
    "L17" = 11.37362843174818
    "N17" = 20.99683581013366
    "B18" = 15.06306221228138
    "D18" = 7.571195650914693
    "E18" = 12.60392322012255
    "F18" = 38.59236697969609
    "G18" = 44.96428115752128
    "H18" = 18.18683317430739
    "I18" = 28.62325312335112
    "J18" = 10.11116699150423
    "K18" = 13.35561950813099
    "L18" = 11.36378437938092
    "N18" = 21.0118406512279
    "B19" = 15.05209478828871
    "D19" = 7.571139698907153
    "E19" = 12.60506145148494
    "F19" = 38.59284018806836
    "G19" = 44.95859909912607
    "H19" = 18.18931703411628
    "I19" = 28.62996535498403
    "J19" = 10.11196134472404
    "K19" = 13.33626009662945
    "L19" = 11.36048047951839
    "N19" = 21.01695110459311
    "B20" = 15.10155622585687
    "D20" = 7.571424732593746
    "E20" = 12.59998447408057
    "F20" = 38.5912683541263
    "G20" = 44.98497180254613
    "H20" = 18.17834453536699
    "I20" = 28.60005012290861
    "J20" = 10.10841109801624
    "K20" = 13.42332321142956
    "L20" = 11.37546413406034
    "N20" = 20.99407302618916
    "B21" = 15.26948560302381
    "D21" = 7.572961924768961
    "E21" = 12.58371587774233
    "F21" = 38.59574121576026
    "G21" = 45.08759037176301
    "H21" = 18.14507056074151
    "I21" = 28.50462532484497
    "J21" = 10.09690940591109
    "K21" = 13.71458143801472
    "L21" = 11.42780507597366
    "N21" = 20.91929237889162
    "B22" = 15.38038818286779
    "D22" = 7.57439573466575
    "E22" = 12.5736699378584
    "F22" = 38.60592571272346
    "G22" = 45.16504500932746
    "H22" = 18.12600355834596
    "I22" = 28.446052449202
    "J22" = 10.08970929355542
    "K22" = 13.90366995399737
    "L22" = 11.46346769069537
    "N22" = 20.87196774309925
    "B23" = 15.32110172762788
    "D23" = 7.573591393663892
    "E23" = 12.57897809181623
    "F23" = 38.59982633946109
    "G23" = 45.12275940686307
    "H23" = 18.13593467750596
    "I23" = 28.47696673651935
    "J23" = 10.09352320088569
    "K23" = 13.80288606541721
    "L23" = 11.44430302303886
    "N23" = 20.89708447597629
    "B24" = 15.09883309136652
    "D24" = 7.571406892817762
    "E24" = 12.60026030329603
    "F24" = 38.5913179605014
    "G24" = 44.98347071050087
    "H24" = 18.17893358506425
    "I24" = 28.60167379438667
    "J24" = 10.10860445276816
    "K24" = 13.4185460220377
    "L24" = 11.37463370193589
    "N24" = 20.9953215150764
    "B25" = 14.86499102632912
    "D25" = 7.570910795486192
    "E25" = 12.6257464701463
    "F25" = 38.61335987528211
    "G25" = 44.87814463745906
    "H25" = 18.23680340219421
    "I25" = 28.75246252133025
    "J25" = 10.12623969962929
    "K25" = 13.00067855967588
    "L25" = 11.30594760753709
    "N25" = 21.10797710480795
}

foreach ($cellRef in $newValues.Keys) {
    $ws.Range($cellRef).Value = $newValues[$cellRef]
}

Write-Output ("Updated " + $newValues.Count + " cells")
